# Timesheet update: "adding 4 listings tot the text"
#
# - Sheet2!C22 formula gains an extra "+4" (9+9+9+9+9+4 -> 9+9+9+9+9+4+4),
#   which ripples through Sheet1's running totals (E23, G23:G29, E31) and
#   the "time left" formula in B32 (which now references F24 instead of
#   F23, matching the day the totals actually shifted on).
# - The commentary shared-string for Sheet2!J22 gets an extra clause about
#   adding listings.
# - Both sheets' saved cursor/selection moves to where the author was
#   last working.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: extend the hour-tally formula for the extra listings work ---
$ws2.Range("C22").Formula = "=9+9+9+9+9+4+4"

# --- Sheet2: update the week's log entry text ---
$ws2.Range("J22").Value = "replacing badly cropt image on front page, rereading guidelines_thesis.pdf, moving caption of tables, add listings to some bug explanations"

# --- Sheet1: B32 now measures remaining time from row 24 onward ---
$ws1.Range("B32").Formula = "=E31-F24"

# --- Restore the saved selections/active cells on each sheet ---
$null = $ws1.Range("B33").Select()
$null = $ws1.Activate()

$null = $ws2.Range("C23").Select()
$null = $ws2.Activate()
